$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# Locate the "IMPORTANT" heading paragraph; the six paragraphs that follow it
# are the ones touched by this edit:
#   1: "Premiere methode peut etre abandonnee."                                 -> (cleared)
#   2: "periode A fixe de 3-6 mois, ..."                                        -> (cleared)
#   3: "" (blank)                                                               -> new note #1
#   4: "" (blank)                                                               -> new note #2
#   5: "se baser sur la recency est fausse vu que la derniere pourrait etre tardive." -> (cleared)
#   6: "" (blank, last paragraph before the section break)                      -> new note #3
$importantIdx = Get-ParaIndexByText $d "IMPORTANT"

$note1 = "- Le probleme avec mes graphs c-est la saisonabilite, tu dois avoir une periode de solde."
$note2 = "C-est pour ca que tu as 2 piques dans ta simulation, ce n-est donc pas vraiment representatif d-un probleme d-obsolescence de ton modele. Ok ?"
$note3 = "Variable categorielle, comment evaluer les distances vu que dire qu-une orange est plus eloignee qu-une pomme est nonsensical. Du coup comme tu as toujours 0 ou 1 ca risque de prendre la priorite sur les autres."

# The three paragraphs that must become text-less again are deleted outright
# (paragraph mark included) rather than having their run text blanked out, so
# the remaining/blank paragraphs keep a normal empty <w:r/> shell instead of
# losing their run entirely. Delete from the last one to the first so the
# earlier indices stay valid while we work.
$p5 = $importantIdx + 5
$p2 = $importantIdx + 2
$p1 = $importantIdx + 1

$rng = $d.Paragraphs.Item($p5).Range
$d.Range($rng.Start, $rng.End).Delete()

$rng = $d.Paragraphs.Item($p2).Range
$d.Range($rng.Start, $rng.End).Delete()

$rng = $d.Paragraphs.Item($p1).Range
$d.Range($rng.Start, $rng.End).Delete()

# After the three deletions above, the paragraphs right after "IMPORTANT" are
# exactly the three that used to be blank (originally in slots 3, 4 and 6).
# Re-insert three fresh blank paragraphs so the overall paragraph count /
# layout is restored to six entries following "IMPORTANT".
$d.Paragraphs.Item($importantIdx + 1).Range.InsertParagraphBefore()
$d.Paragraphs.Item($importantIdx + 1).Range.InsertParagraphBefore()
$d.Paragraphs.Item($importantIdx + 5).Range.InsertParagraphBefore()

# Fill the three originally-blank paragraphs (now at offsets 3, 4 and 6) with
# the new note text.
$d.Paragraphs.Item($importantIdx + 3).Range.Text = $note1
$d.Paragraphs.Item($importantIdx + 4).Range.Text = $note2
$d.Paragraphs.Item($importantIdx + 6).Range.Text = $note3
